$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain text (e.g. "37.830.18", "0.626").
# Excel COM auto-converts plain numeric-looking text assigned to
# .Value into a real Number, so pin those specific cells to Text
# format first (only the ones whose new value is unambiguous, to avoid
# touching the format of cells that are not otherwise edited).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "37.823.41"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.091.92"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "233.86"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "0.625"
$ws.Range("D7").Value = "58.54"
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("D12").Value = "15.29"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").Value = "2.400.82"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "21.28"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "2.090.01"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "37.795.20"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "71.07"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "229.96"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "9.83"
$ws.Range("E26").Value = "  +8.95%  "
$ws.Range("D27").Value = "171.38"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").Value = "19.52"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "4.71"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "0.0633"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "3.34"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("E40").Value = "  +9.42%  "
$ws.Range("D41").Value = "101.44"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "1.18"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "16.79"
$ws.Range("E45").Value = "  +5.00%  "
$ws.Range("D46").Value = "1.452.67"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "4.14"
$ws.Range("E47").Value = "  -4.50%  "
$ws.Range("D48").Value = "1.06"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "7.25"
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").Value = "2.98"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").Value = "2.284.50"
$ws.Range("E51").Value = "  +0.21%  "
